$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.966.64"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.307.62"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.299.82"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.627"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "3.834.99"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "3.309.27"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "63.938.61"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.978"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.93%  "
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "60.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "561.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").Value = "0.0₃0722"
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").Value = "3.034.68"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.132"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.22%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("E51").Value = "  -1.05%  "
